# Shard edit: add a "State" column to hotel_info (between Hotel_Name and City)
# and move review_info ahead of hotel_info in sheet order.

$wb = $excel.ActiveWorkbook

$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")

# Insert a new column C (State) right before the existing City column,
# shifting City/Zip/... one column to the right.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"

# review_info now comes before hotel_info in the workbook's tab order.
$wsReview.Move($wsHotel)
